$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-07 Thursday" "2024-03-08 Friday"

Replace-Text "124×9=1116" "470×9=4230"
Replace-Text "509×7=3563" "676×8=5408"
Replace-Text "160×3=480" "786×5=3930"
Replace-Text "330×8=2640" "290×6=1740"
Replace-Text "639×2=1278" "649×8=5192"

Replace-Text "888×9=7992" "156×5=780"
Replace-Text "426×8=3408" "491×6=2946"
Replace-Text "854×6=5124" "795×5=3975"
Replace-Text "882×9=7938" "992×5=4960"
Replace-Text "598×9=5382" "290×8=2320"

Replace-Text "263×4=1052" "498×5=2490"
Replace-Text "648×9=5832" "670×3=2010"
Replace-Text "224×4=896" "730×7=5110"
Replace-Text "549×2=1098" "467×2=934"
Replace-Text "370×7=2590" "206×8=1648"

Replace-Text "934×6=5604" "546×5=2730"
Replace-Text "541×9=4869" "748×9=6732"
Replace-Text "537×8=4296" "972×6=5832"
Replace-Text "673×3=2019" "232×5=1160"
Replace-Text "199×7=1393" "869×3=2607"

Replace-Text "716×9=6444" "120×3=360"
Replace-Text "149×4=596" "300×6=1800"
Replace-Text "383×6=2298" "633×7=4431"
Replace-Text "487×7=3409" "920×6=5520"
Replace-Text "514×6=3084" "960×7=6720"

Write-Output "done"
